$d = $word.ActiveDocument

# Ordered list of (old, new) text replacements for the table cells.
# NOTE: the pair that produces "64÷5=12, 4" must run AFTER the pair that
# consumes "64÷5=12, 4" as its search text, to avoid a false re-match.
$pairs = @(
    @("64÷5=12, 4", "32÷7=4, 4"),
    @("16÷6=2, 4", "57÷7=8, 1"),
    @("20÷4=5, 0", "48÷3=16, 0"),
    @("29÷8=3, 5", "10÷6=1, 4"),
    @("47÷8=5, 7", "72÷8=9, 0"),
    @("50÷8=6, 2", "45÷8=5, 5"),
    @("88÷7=12, 4", "33÷4=8, 1"),
    @("20÷7=2, 6", "85÷3=28, 1"),
    @("32÷3=10, 2", "13÷6=2, 1"),
    @("84÷3=28, 0", "72÷6=12, 0"),
    @("48÷9=5, 3", "64÷5=12, 4"),
    @("86÷9=9, 5", "77÷8=9, 5"),
    @("82÷6=13, 4", "75÷5=15, 0"),
    @("21÷4=5, 1", "74÷2=37, 0"),
    @("61÷8=7, 5", "47÷2=23, 1"),
    @("71÷9=7, 8", "97÷8=12, 1"),
    @("42÷6=7, 0", "91÷2=45, 1"),
    @("53÷8=6, 5", "70÷3=23, 1"),
    @("76÷8=9, 4", "55÷7=7, 6"),
    @("78÷3=26, 0", "91÷7=13, 0"),
    @("57÷2=28, 1", "76÷6=12, 4"),
    @("45÷5=9, 0", "11÷6=1, 5"),
    @("48÷8=6, 0", "91÷6=15, 1"),
    @("75÷3=25, 0", "54÷3=18, 0"),
    @("40÷7=5, 5", "90÷3=30, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
